$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 3356
    $ws.Range("F4").Value = 90
    $ws.Range("F5").Value = 669
}
